{"js": "// Update the worksheet date and the 25 two-digit x two-digit multiplication\n// problems/answers in the table to a newly generated set.\nconst replacements = [\n  { find: \"2025-03-10 Monday\", replace: \"2025-03-11 Tuesday\" },\n  { find: \"52\u00d777=4004\", replace: \"96\u00d740=3840\" },\n  { find: \"68\u00d738=2584\", replace: \"90\u00d743=3870\" },\n  { find: \"93\u00d756=5208\", replace: \"34\u00d779=2686\" },\n  { find: \"95\u00d716=1520\", replace: \"34\u00d782=2788\" },\n  { find: \"13\u00d748=624\", replace: \"83\u00d764=5312\" },\n  { find: \"87\u00d776=6612\", replace: \"38\u00d764=2432\" },\n  { find: \"69\u00d797=6693\", replace: \"39\u00d764=2496\" },\n  { find: \"63\u00d786=5418\", replace: \"99\u00d728=2772\" },\n  { find: \"84\u00d788=7392\", replace: \"13\u00d763=819\" },\n  { find: \"21\u00d720=420\", replace: \"13\u00d712=156\" },\n  { find: \"95\u00d748=4560\", replace: \"67\u00d732=2144\" },\n  { find: \"72\u00d794=6768\", replace: \"30\u00d725=750\" },\n  { find: \"68\u00d776=5168\", replace: \"66\u00d774=4884\" },\n  { find: \"68\u00d725=1700\", replace: \"30\u00d717=510\" },\n  { find: \"11\u00d789=979\", replace: \"95\u00d721=1995\" },\n  { find: \"81\u00d761=4941\", replace: \"71\u00d714=994\" },\n  { find: \"63\u00d744=2772\", replace: \"60\u00d786=5160\" },\n  { find: \"93\u00d761=5673\", replace: \"16\u00d791=1456\" },\n  { find: \"79\u00d712=948\", replace: \"80\u00d729=2320\" },\n  { find: \"61\u00d719=1159\", replace: \"87\u00d793=8091\" },\n  { find: \"18\u00d742=756\", replace: \"46\u00d767=3082\" },\n  { find: \"25\u00d786=2150\", replace: \"84\u00d767=5628\" },\n  { find: \"27\u00d751=1377\", replace: \"55\u00d746=2530\" },\n  { find: \"55\u00d731=1705\", replace: \"86\u00d781=6966\" },\n  { find: \"61\u00d771=4331\", replace: \"51\u00d799=5049\" },\n];\n\nconst body = context.document.body;\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the 25 two-digit x two-digit multiplication\n# problems/answers in the table to a newly generated set.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Find = '2025-03-10 Monday'; Replace = '2025-03-11 Tuesday'},\n    @{Find = '52\u00d777=4004'; Replace = '96\u00d740=3840'},\n    @{Find = '68\u00d738=2584'; Replace = '90\u00d743=3870'},\n    @{Find = '93\u00d756=5208'; Replace = '34\u00d779=2686'},\n    @{Find = '95\u00d716=1520'; Replace = '34\u00d782=2788'},\n    @{Find = '13\u00d748=624'; Replace = '83\u00d764=5312'},\n    @{Find = '87\u00d776=6612'; Replace = '38\u00d764=2432'},\n    @{Find = '69\u00d797=6693'; Replace = '39\u00d764=2496'},\n    @{Find = '63\u00d786=5418'; Replace = '99\u00d728=2772'},\n    @{Find = '84\u00d788=7392'; Replace = '13\u00d763=819'},\n    @{Find = '21\u00d720=420'; Replace = '13\u00d712=156'},\n    @{Find = '95\u00d748=4560'; Replace = '67\u00d732=2144'},\n    @{Find = '72\u00d794=6768'; Replace = '30\u00d725=750'},\n    @{Find = '68\u00d776=5168'; Replace = '66\u00d774=4884'},\n    @{Find = '68\u00d725=1700'; Replace = '30\u00d717=510'},\n    @{Find = '11\u00d789=979'; Replace = '95\u00d721=1995'},\n    @{Find = '81\u00d761=4941'; Replace = '71\u00d714=994'},\n    @{Find = '63\u00d744=2772'; Replace = '60\u00d786=5160'},\n    @{Find = '93\u00d761=5673'; Replace = '16\u00d791=1456'},\n    @{Find = '79\u00d712=948'; Replace = '80\u00d729=2320'},\n    @{Find = '61\u00d719=1159'; Replace = '87\u00d793=8091'},\n    @{Find = '18\u00d742=756'; Replace = '46\u00d767=3082'},\n    @{Find = '25\u00d786=2150'; Replace = '84\u00d767=5628'},\n    @{Find = '27\u00d751=1377'; Replace = '55\u00d746=2530'},\n    @{Find = '55\u00d731=1705'; Replace = '86\u00d781=6966'},\n    @{Find = '61\u00d771=4331'; Replace = '51\u00d799=5049'},\n)\n\nforeach ($r in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    # Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, MatchSoundsLike,\n    #         MatchAllWordForms, Forward, Wrap, Format, ReplaceWith, Replace)\n    # Wrap:=1 (wdFindContinue), Replace:=2 (wdReplaceAll)\n    $find.Execute($r.Find, $true, $false, $false, $false, $false, $true, 1, $false, $r.Replace, 2)\n}\n"}
